$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(class_weight=''balanced'',
                                                                     criterion=''entropy'',
                                                                     max_depth=5,
                                                                     max_features=''sqrt'',
                                                                     min_samples_leaf=6,
                                                                     min_samples_split=6,
                                                                     random_state=42),
                                    random_state=42))])'
$ws.Range("B2").Value = 0.5416028416028416
$ws.Range("C2").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': RobustScaler(), ''model__n_estimators'': 50, ''model__estimator__min_samples_split'': 6, ''model__estimator__min_samples_leaf'': 6, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 5, ''model__estimator__criterion'': ''entropy'', ''model__estimator__class_weight'': ''balanced''}'
$ws.Range("D2").Value = 0.7999999999999999
$ws.Range("E2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1]'
$ws.Range("F2").Value = '[1 0 1 1 1 1 0 0 1 0 0 1]'
$ws.Range("G2").Value = 42
$ws.Range("H2").Value = 0.9006177907956395
$ws.Range("I2").Value = 0.02356486820354882
$ws.Range("J2").Value = 0.4554321123321124
$ws.Range("K2").Value = 0.1668581131322856

# Row 3
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(max_depth=4,
                                                                     max_features=''sqrt'',
                                                                     min_samples_leaf=4,
                                                                     min_samples_split=5,
                                                                     random_state=42),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B3").Value = 0.631068931068931
$ws.Range("C3").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': RobustScaler(), ''model__n_estimators'': 10, ''model__estimator__min_samples_split'': 5, ''model__estimator__min_samples_leaf'': 4, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 4, ''model__estimator__criterion'': ''gini'', ''model__estimator__class_weight'': None}'
$ws.Range("D3").Value = 0.5714285714285714
$ws.Range("E3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0]'
$ws.Range("F3").Value = '[0 1 0 1 1 1 0 1 0 1 1 0]'
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.8908199617682057
$ws.Range("I3").Value = 0.02548209242292365
$ws.Range("J3").Value = 0.5392226218226218
$ws.Range("K3").Value = 0.1565767270107846

# Row 4
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(class_weight=''balanced'',
                                                                     max_depth=6,
                                                                     min_samples_leaf=6,
                                                                     min_samples_split=6,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B4").Value = 0.6902219349278171
$ws.Range("C4").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__min_samples_split'': 6, ''model__estimator__min_samples_leaf'': 6, ''model__estimator__max_features'': None, ''model__estimator__max_depth'': 6, ''model__estimator__criterion'': ''gini'', ''model__estimator__class_weight'': ''balanced''}'
$ws.Range("D4").Value = 0.4285714285714285
$ws.Range("E4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0]'
$ws.Range("F4").Value = '[0 1 1 1 0 0 0 1 1 0 0 1]'
$ws.Range("G4").Value = 23
$ws.Range("H4").Value = 0.910592479779327
$ws.Range("I4").Value = 0.02195213284173735
$ws.Range("J4").Value = 0.5409538262391203
$ws.Range("K4").Value = 0.1475826440591196

# Row 5
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(class_weight=''balanced'',
                                                                     max_depth=6,
                                                                     max_features=''sqrt'',
                                                                     min_samples_leaf=3,
                                                                     min_samples_split=4,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B5").Value = 0.6918181818181818
$ws.Range("C5").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__min_samples_split'': 4, ''model__estimator__min_samples_leaf'': 3, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 6, ''model__estimator__criterion'': ''gini'', ''model__estimator__class_weight'': ''balanced''}'
$ws.Range("D5").Value = 0.8000000000000002
$ws.Range("E5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1]'
$ws.Range("F5").Value = '[0 1 1 0 0 1 0 0 0 1 1 0]'
$ws.Range("G5").Value = 99
$ws.Range("H5").Value = 0.8738143952987326
$ws.Range("I5").Value = 0.02355277637360543
$ws.Range("J5").Value = 0.5291931382343147
$ws.Range("K5").Value = 0.1640374644108655
